# Loan RBI, Variable Instalments
# The "Repayment schedule" sheet gets a new (blank) column inserted before
# the existing "Late" column (column N), pushing the old N/O/P ("Late",
# "heading"/Outstanding-label, "Outstanding") columns one slot to the right
# (O/P/Q). The sheet also becomes the active tab/selection, replacing
# "NewLoanInput" as the selected sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at N (shifts old N->O, O->P, P->Q).
$ws.Columns("N:N").Insert()

# The newly inserted column picks up the width of its left neighbour (M),
# matching Excel's normal "insert column" behaviour.
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab, with a single cell (R8)
# selected (this also clears "NewLoanInput"'s previous tabSelected flag).
$ws.Activate()
$ws.Range("R8").Select() | Out-Null
